# Add 6 new registration rows to the data table, bringing it from
# A1:C35 to A1:C41. Two of the new records land above the existing
# "+553195267711" row (pushing it from row 2 to row 4) and four more
# land directly below it (pushing the rest of the table, previously
# starting at old row 3, down to row 9 and on).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows right after the header (old row 2 -> row 4).
$ws.Rows("2:3").Insert()

# Insert 4 more blank rows right after the relocated row (old row 3
# -> row 9, and everything below follows).
$ws.Rows("5:8").Insert()

# The inserted rows pick up neighboring formatting inconsistently (some
# inherit the bold/red header look), so restore the standard data-row
# style (the one used throughout the table, e.g. row 9) across the
# whole new block.
$ws.Range("A9:C9").Copy()
$ws.Range("A2:C3").PasteSpecial(-4122)
$ws.Range("A5:C8").PasteSpecial(-4122)

# New records (Telefone, DDD, Data Inscricao). A leading apostrophe forces
# Excel to keep these digit-heavy values as literal text instead of
# auto-converting them to numbers/dates.
$newRows = @{
    2 = @("+553174012843", "31", "2024-12-12")
    3 = @("+555496738305", "54", "2024-12-10")
    5 = @("+555180524794", "51", "2024-11-14")
    6 = @("+5519997335593", "19", "2024-11-11")
    7 = @("+555391771383", "53", "2024-11-10")
    8 = @("+553173576911", "31", "2024-11-08")
}

foreach ($r in $newRows.Keys) {
    $vals = $newRows[$r]
    $ws.Cells.Item($r, 1).Value = "'" + $vals[0]
    $ws.Cells.Item($r, 2).Value = "'" + $vals[1]
    $ws.Cells.Item($r, 3).Value = "'" + $vals[2]
}
